$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.608105301856995
$ws.Range("B1").Value = 2.790774822235107
$ws.Range("C1").Value = 4.983926296234131
$ws.Range("D1").Value = 1.461317539215088
$ws.Range("E1").Value = 0.8361047506332397
